# Auto-generated edit script applying the Marilith_Profits.xlsx cell-value diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 261.16666
$ws.Range("I12").Value = 288.4
$ws.Range("J12").Value = 125
$ws.Range("K12").Value = 288.4
$ws.Range("L12").Value = 125
$ws.Range("M12").Value = -118.4
$ws.Range("N12").Value = -465
$ws.Range("H33").Value = 953.625
$ws.Range("I33").Value = 96.61539
$ws.Range("J33").Value = 4667.3335
$ws.Range("K33").Value = 96.61539
$ws.Range("L33").Value = 4667.3335
$ws.Range("M33").Value = 132.38461
$ws.Range("N33").Value = -5125.3335
$ws.Range("H62").Value = 5142.9287
$ws.Range("J62").Value = 5999.2
$ws.Range("L62").Value = 5999.2
$ws.Range("N62").Value = -7247.2
$ws.Range("H65").Value = 5142.9287
$ws.Range("J65").Value = 5999.2
$ws.Range("L65").Value = 29996
$ws.Range("N65").Value = -36236
$ws.Range("H98").Value = 3644.3333
$ws.Range("I98").Value = 3599.875
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 3599.875
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -2101.875
$ws.Range("N98").Value = -6996
$ws.Range("H100").Value = 858
$ws.Range("I100").Value = 862.25
$ws.Range("J100").Value = 849.5
$ws.Range("K100").Value = 862.25
$ws.Range("L100").Value = 849.5
$ws.Range("M100").Value = -321.25
$ws.Range("N100").Value = -1931.5
$ws.Range("H114").Value = 79999
$ws.Range("J114").Value = 79999
$ws.Range("L114").Value = 79999
$ws.Range("N114").Value = -88677
$ws.Range("H122").Value = 3644.3333
$ws.Range("I122").Value = 3599.875
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10799.625
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8349.625
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 6788
$ws.Range("I132").Value = 6987.5
$ws.Range("K132").Value = 20962.5
$ws.Range("M132").Value = -18432.5
$ws.Range("H138").Value = 3051
$ws.Range("J138").Value = 3771.25
$ws.Range("L138").Value = 11313.75
$ws.Range("N138").Value = -21593.75
$ws.Range("H141").Value = 5359.8887
$ws.Range("I141").Value = 5087.1177
$ws.Range("K141").Value = 15261.3531
$ws.Range("M141").Value = -10081.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 237.5
$ws.Range("I2").Value = 215
$ws.Range("K2").Value = 215
$ws.Range("M2").Value = -102
$ws.Range("H32").Value = 7475
$ws.Range("I32").Value = 7475
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7475
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7188
$ws.Range("N32").ClearContents()
$ws.Range("H116").Value = 237.5
$ws.Range("I116").Value = 215
$ws.Range("K116").Value = 215
$ws.Range("M116").Value = 2079
$ws.Range("H132").Value = 2329.6667
$ws.Range("I132").Value = 2095.65
$ws.Range("K132").Value = 6286.950000000001
$ws.Range("M132").Value = -3756.950000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 237.5
$ws.Range("I3").Value = 215
$ws.Range("K3").Value = 215
$ws.Range("M3").Value = -101
$ws.Range("H26").Value = 32000
$ws.Range("I26").Value = 32000
$ws.Range("K26").Value = 32000
$ws.Range("M26").Value = -31708
$ws.Range("H86").Value = 3379.7
$ws.Range("I86").Value = 3422
$ws.Range("J86").Value = 2999
$ws.Range("K86").Value = 3422
$ws.Range("L86").Value = 2999
$ws.Range("M86").Value = -2299
$ws.Range("N86").Value = -5245
$ws.Range("H89").Value = 3379.7
$ws.Range("I89").Value = 3422
$ws.Range("J89").Value = 2999
$ws.Range("K89").Value = 17110
$ws.Range("L89").Value = 14995
$ws.Range("M89").Value = -11494
$ws.Range("N89").Value = -26227
$ws.Range("H94").Value = 1135.7084
$ws.Range("I94").Value = 1132.4783
$ws.Range("J94").Value = 1210
$ws.Range("K94").Value = 1132.4783
$ws.Range("L94").Value = 1210
$ws.Range("M94").Value = -681.4783
$ws.Range("N94").Value = -2112
$ws.Range("H96").Value = 20249.75
$ws.Range("I96").Value = 20249.75
$ws.Range("K96").Value = 20249.75
$ws.Range("M96").Value = -17503.75
$ws.Range("H107").Value = 653.46155
$ws.Range("I107").Value = 653.46155
$ws.Range("K107").Value = 653.46155
$ws.Range("M107").Value = 1266.53845
$ws.Range("H134").Value = 6470.75
$ws.Range("I134").Value = 7144.9
$ws.Range("K134").Value = 21434.7
$ws.Range("M134").Value = -18899.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 373.22223
$ws.Range("I22").Value = 100
$ws.Range("J22").Value = 509.83334
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = 509.83334
$ws.Range("M22").Value = 250
$ws.Range("N22").Value = -1209.83334
$ws.Range("H105").Value = 951.8889
$ws.Range("I105").Value = 951.8889
$ws.Range("K105").Value = 951.8889
$ws.Range("M105").Value = 795.1111

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 97938
$ws.Range("J37").Value = 97938
$ws.Range("L37").Value = 293814
$ws.Range("N37").Value = -294038
$ws.Range("H81").Value = 7562.1113
$ws.Range("J81").Value = 8257.375
$ws.Range("L81").Value = 24772.125
$ws.Range("N81").Value = -27018.125
$ws.Range("H84").Value = 7562.1113
$ws.Range("J84").Value = 8257.375
$ws.Range("L84").Value = 74316.375
$ws.Range("N84").Value = -85548.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 106665
$ws.Range("J69").Value = 106665
$ws.Range("L69").Value = 106665
$ws.Range("N69").Value = -108163
$ws.Range("H70").Value = 2763.3333
$ws.Range("I70").Value = 2745
$ws.Range("J70").Value = 2800
$ws.Range("K70").Value = 2745
$ws.Range("L70").Value = 2800
$ws.Range("M70").Value = -2475
$ws.Range("N70").Value = -3340
$ws.Range("H72").Value = 106665
$ws.Range("J72").Value = 106665
$ws.Range("L72").Value = 319995
$ws.Range("N72").Value = -327483
$ws.Range("H73").Value = 2763.3333
$ws.Range("I73").Value = 2745
$ws.Range("J73").Value = 2800
$ws.Range("K73").Value = 2745
$ws.Range("L73").Value = 2800
$ws.Range("M73").Value = -1809
$ws.Range("N73").Value = -4672
$ws.Range("H80").Value = 5218.75
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 5218.75
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H97").Value = 627.4400000000001
$ws.Range("I97").Value = 660.7143
$ws.Range("K97").Value = 660.7143
$ws.Range("M97").Value = -164.7143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10845.2
$ws.Range("I132").Value = 10630.134
$ws.Range("J132").Value = 11490.4
$ws.Range("K132").Value = 31890.402
$ws.Range("L132").Value = 34471.2
$ws.Range("M132").Value = -29360.402
$ws.Range("N132").Value = -39531.2
$ws.Range("H136").Value = 3639.3635
$ws.Range("J136").Value = 4068.3333
$ws.Range("L136").Value = 12204.9999
$ws.Range("N136").Value = -17304.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 899.5
$ws.Range("I100").Value = 800.5
$ws.Range("K100").Value = 1601
$ws.Range("M100").Value = -1060
$ws.Range("H122").Value = 615.06665
$ws.Range("I122").Value = 615.06665
$ws.Range("K122").Value = 1845.19995
$ws.Range("M122").Value = 604.8000500000001
$ws.Range("H130").Value = 80249.875
$ws.Range("J130").Value = 80249.875
$ws.Range("L130").Value = 80249.875
$ws.Range("N130").Value = -90289.875
$ws.Range("H132").Value = 768.8
$ws.Range("I132").Value = 662.25
$ws.Range("K132").Value = 1986.75
$ws.Range("M132").Value = 543.25
